# testUserDoc2.docx : add "Debut "/" fin" markers around the userdoc
# field pair and drop the "id=" AQL parameter from the m:userdoc field
# code (keep only the 'value1' expression), per the M2Doc "[REA] Add
# support of user part document" change.
#
# The paragraph mixes plain w:r/w:t runs with w:fldChar / w:instrText
# runs (two fields: "m:userdoc ..." ... "m:enduserdoc") plus a
# _GoBack bookmark. Find/Replace does not see inside field codes, and
# reordering the runs (the "end" fldChar + result text + second field
# now come BEFORE the bookmark instead of after) is not expressible as
# a simple text substitution, so the whole paragraph content is
# replaced in one shot via Range.InsertXML with the exact OOXML we
# want - this keeps the two fields, the bookmark and the field-code
# text all correctly recognised afterwards.

$d = $word.ActiveDocument

# Locate the paragraph that hosts the "m:userdoc" field (rather than
# hard-coding an index) so the script is resilient to unrelated
# paragraphs being added/removed earlier in the document.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Fields.Count -gt 0) {
        $f = $p.Range.Fields.Item(1)
        if ($f.Code.Text -like "*userdoc*") {
            $target = $p
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate the userdoc field paragraph"
}

$newParagraphXml = @'
<w:p w:rsidR="005C73CF" w:rsidRDefault="005C73CF" w:rsidP="005C73CF"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Début </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="0056766F"><w:instrText>m</w:instrText></w:r><w:r w:rsidR="00A7781B"><w:instrText>:userdoc</w:instrText></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:instrText xml:space="preserve">'value1' </w:instrText></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidR="00A7781B"><w:t>User document part Texte 1</w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="0056766F"><w:instrText>m</w:instrText></w:r><w:r w:rsidR="00A7781B"><w:instrText>:enduserdoc</w:instrText></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:r><w:t xml:space="preserve"> fin</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newParagraphXml + '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($packageXml)
